# "method for deleting technicians are expanded"
# Adds a rowNumber column and a positiveMessage column to the technician
# test data on the UserManagement sheet, and fills in the previously-blank
# SI_003 row with a third positive test case.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UserManagement")

# --- 1) Insert the two new columns -----------------------------------
# First a column before the technician block (becomes "rowNumber"),
# then - after that shift - a column before the empty-message block
# (becomes "positiveMessage").
$ws.Range("K1").EntireColumn.Insert()
$ws.Range("P1").EntireColumn.Insert()

# --- 2) Header row (row 3) labels for the two new columns ------------
$ws.Range("P3").Value2 = "positiveMessage"
$ws.Range("K3").Value2 = "rowNumber"

# --- 3) Row 4 (SI_001): rowNumber / positiveMessage values ------------
$ws.Range("P4").Value2 = "successfully created."
$ws.Range("K4").Value2 = "4"
$ws.Range("L4").Value2 = "1"

# --- 4) Row 5 (SI_002): rowNumber / positiveMessage values -------------
$ws.Range("K5").Value2 = "4"
$ws.Range("P5").Value2 = "successfully created."

# --- 5) Row 6 (SI_003) was essentially empty before; fill it in with --
#        the same data shape as row 4 (a third positive user-management
#        test case), keeping A6's existing value/style untouched.
$ws.Range("B6").Value2 = $ws.Range("B4").Value2
$ws.Range("C6").Value2 = $ws.Range("C4").Value2
$ws.Range("D6").Value2 = $ws.Range("D4").Value2
$ws.Range("E6").Value2 = $ws.Range("E4").Value2
$ws.Range("F6").Value2 = $ws.Range("F4").Value2
$ws.Range("G6").Value2 = $ws.Range("G4").Value2
$ws.Range("I6").Value2 = $ws.Range("I4").Value2
$ws.Range("J6").Value2 = $ws.Range("J4").Value2
$ws.Range("K6").Value2 = $ws.Range("K4").Value2
$ws.Range("L6").Value2 = "3"
$ws.Range("P6").Value2 = $ws.Range("P4").Value2
$ws.Range("Q6").Value2 = $ws.Range("Q4").Value2
$ws.Range("R6").Value2 = $ws.Range("R4").Value2
$ws.Range("S6").Value2 = $ws.Range("S4").Value2
$ws.Range("U6").Value2 = $ws.Range("U4").Value2

# --- 6) Row 7 (SI_004) also gets the positiveMessage value -------------
$ws.Range("P7").Value2 = "successfully created."

# --- 7) Restyle the technician-number cells: writing a new value onto --
#        a "quote-prefixed text" cell drops that formatting, so restore
#        it afterwards by copying the format from an unmodified sibling
#        cell that still carries the correct style.
$ws.Range("M4").Copy()
$ws.Range("K4").PasteSpecial(-4122)
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("K5").PasteSpecial(-4122)
$ws.Range("K7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 6 (B:U) picks up row 4's formatting pattern in one shot.
$ws.Range("B4:U4").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 8) Mirror the edited session's selection --------------------------
$ws.Range("K14").Select()
